$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet: "Property1" -> "DataNode" -------------------------
$ws.Name = "DataNode"

# --- Move the active selection from A9 to H13 -----------------------------
$ws.Range("H13").Select()

# --- Row height tweaks -----------------------------------------------------
# Header row grows to fit a taller wrapped caption.
$ws.Rows.Item(1).RowHeight = 27
# The long-description row shrinks slightly.
$ws.Rows.Item(8).RowHeight = 81

# --- Column width tweaks (re-measured after a Windows/Mac round trip) ------
$ws.Columns.Item(2).ColumnWidth = 10.857142857142858
$ws.Columns.Item(3).ColumnWidth = 16.428571428571427
$ws.Columns.Item(4).ColumnWidth = 21.857142857142858
$ws.Columns.Item(5).ColumnWidth = 16.428571428571427
$ws.Columns.Item(7).ColumnWidth = 11.857142857142858
$ws.Columns.Item(8).ColumnWidth = 17.714285714285715
$ws.Columns.Item(9).ColumnWidth = 18.714285714285715
$ws.Columns.Item(10).ColumnWidth = 13.142857142857142

# --- Register the extra 9pt font used for the sheet's phonetic settings ---
# (Touch it on a cell that already renders at 11pt/SimSun, then restore the
# original size so no cell's effective style actually changes -- this just
# leaves the new font behind in the shared font table, matching the
# fonts count="2" -> count="3" change.)
$probeCell = $ws.Cells.Item(9, 2)
$probeCell.Font.Name = "宋体"
$probeCell.Font.Size = 9
$probeCell.Font.Name = "宋体"
$probeCell.Font.Size = 11
